# Calendar session page cleanup: rows 4/5 and 10/11 swap their course
# info (ELT/FONDAMENTI DI ROBOTICA <-> ATM/SISTEMI INFORMATICI), and a
# handful of "Date" values (column I) are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Course Code) and I (Date) hold values that look numeric /
# date-like but must stay stored as plain text, matching the rest of
# the sheet. Force text formatting on those specific cells before
# writing so Excel doesn't reinterpret them as a number / serial date.
$textCells = @("B4","B5","B10","B11","I2","I3","I7","I10","I11","I12","I13")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2: date correction ---
$ws.Range("I2").Value = "2023-07-05"

# --- Row 3: date correction ---
$ws.Range("I3").Value = "2023-07-10"

# --- Row 4: now the ATM / SISTEMI INFORMATICI session ---
$ws.Range("A4").Value = "ATM"
$ws.Range("B4").Value = "85743"
$ws.Range("C4").Value = "SISTEMI INFORMATICI"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = "Gatti Nicola-Mottola Luca"

# --- Row 5: now the ATM / SISTEMI INFORMATICI session ---
$ws.Range("A5").Value = "ATM"
$ws.Range("B5").Value = "85743"
$ws.Range("C5").Value = "SISTEMI INFORMATICI"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = "Gatti Nicola-Mottola Luca"

# --- Row 7: date correction ---
$ws.Range("I7").Value = "2023-07-13"

# --- Row 10: now the ELT / FONDAMENTI DI ROBOTICA session ---
$ws.Range("A10").Value = "ELT"
$ws.Range("B10").Value = "85754"
$ws.Range("C10").Value = "FONDAMENTI DI ROBOTICA"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = "Zanchettin Andrea Maria-Rocco Paolo"
$ws.Range("I10").Value = "2023-07-03"

# --- Row 11: now the ELT / FONDAMENTI DI ROBOTICA session ---
$ws.Range("A11").Value = "ELT"
$ws.Range("B11").Value = "85754"
$ws.Range("C11").Value = "FONDAMENTI DI ROBOTICA"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = "Zanchettin Andrea Maria-Rocco Paolo"
$ws.Range("I11").Value = "2023-07-08"

# --- Row 12: date correction ---
$ws.Range("I12").Value = "2023-07-04"

# --- Row 13: date correction ---
$ws.Range("I13").Value = "2023-07-12"
